$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the INTRON/EXON row (row 8) with the new "Highlighting Exon Task" observation.
$ws.Range("B8").Value = "Does not need a function"
$ws.Range("C8").Value = "Does not need a function."

# Minor wording tweak on the CODING_DNA row.
$ws.Range("D9").Value = "The coding sequence; concatenated exons (1 string)."

$ws.Range("E8").Value = "This task can be completed by only giving the front end the full DNA sequence and an hash of exons positions/length; both retrievable with queries; all the front end needs is to know which sbstring within the main string to highligh; same we decided to do for the restriction sites task."
$ws.Range("D8").Value = "The front end can do this task with 2 lines of code by only having the full DNA sequence and an hash of exons position. We aleady have SQL fnctions retrieving those for the front end to use.  The front end just need to place tags around each exon subtring."

# F8 gets a new "Done" status - copy formatting from the neighboring Done cell (F9) first.
$ws.Range("F9").Copy()
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("F8").Value = "Done"

# Row 8 grows taller to fit the new, longer notes.
$ws.Rows.Item(8).RowHeight = 75

# Update the view: scrolled up/left and a new active selection.
$ws.Range("D10").Select()
